# Update workbook to reflect newer COVID-19 data snapshot (14 Aug 2020, 19:14)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp string (row 1)
$ws.Range("A1").Value = "Datos actualizados a 14 de Agosto de 2020 a las 19:14"

# Update country rows: name (col A) + statistics (cols B:H -> Casos totales, Nuevos casos,
# Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)

$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 5433129
$ws.Range("C4").Value = 17463
$ws.Range("D4").Value = 2845717
$ws.Range("E4").Value = 2416595
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 402
$ws.Range("H4").Value = 170817

$ws.Range("A13").Value = "España"
$ws.Range("B13").Value = 358843
$ws.Range("C13").Value = 2987
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 12
$ws.Range("H13").Value = 28617

$ws.Range("A21").Value = "Turquia"
$ws.Range("B21").Value = 246861
$ws.Range("C21").Value = 1226
$ws.Range("D21").Value = 228980
$ws.Range("E21").Value = 11947
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 22
$ws.Range("H21").Value = 5934

$ws.Range("A24").Value = "Irak"
$ws.Range("B24").Value = 168290
$ws.Range("C24").Value = 4013
$ws.Range("D24").Value = 120129
$ws.Range("E24").Value = 42452
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 68
$ws.Range("H24").Value = 5709

$ws.Range("A27").Value = "Canada"
$ws.Range("B27").Value = 121414
$ws.Range("C27").Value = 180
$ws.Range("D27").Value = 107793
$ws.Range("E27").Value = 4602
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 9019

$ws.Range("A36").Value = "Republica Dominicana"
$ws.Range("B36").Value = 84488
$ws.Range("C36").Value = 1354
$ws.Range("D36").Value = 49539
$ws.Range("E36").Value = 33540
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 16
$ws.Range("H36").Value = 1409

$ws.Range("A37").Value = "Suecia"
$ws.Range("B37").Value = 84294
$ws.Range("C37").Value = 0
$ws.Range("D37").Value = 0
$ws.Range("E37").Value = 0
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 7
$ws.Range("H37").Value = 5783

$ws.Range("A59").Value = "Argelia"
$ws.Range("B59").Value = 37664
$ws.Range("C59").Value = 477
$ws.Range("D59").Value = 26308
$ws.Range("E59").Value = 10005
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = 10
$ws.Range("H59").Value = 1351

$ws.Range("A60").Value = "Afganistan"
$ws.Range("B60").Value = 37431
$ws.Range("C60").Value = 7
$ws.Range("D60").Value = 26714
$ws.Range("E60").Value = 9354
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 1363

$ws.Range("A67").Value = "Etiopia"
$ws.Range("B67").Value = 27242
$ws.Range("C67").Value = 1038
$ws.Range("D67").Value = 11660
$ws.Range("E67").Value = 15090
$ws.Range("F67").Value = 0
$ws.Range("G67").Value = 13
$ws.Range("H67").Value = 492

$ws.Range("A68").Value = "Irlanda"
$ws.Range("B68").Value = 26995
$ws.Range("C68").Value = 66
$ws.Range("D68").Value = 23364
$ws.Range("E68").Value = 1857
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 1774

$ws.Range("A94").Value = "Libano"
$ws.Range("B94").Value = 8045
$ws.Range("C94").Value = 334
$ws.Range("D94").Value = 2551
$ws.Range("E94").Value = 5400
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = 2
$ws.Range("H94").Value = 94

$ws.Range("A95").Value = "Tayikistan"
$ws.Range("B95").Value = 7950
$ws.Range("C95").Value = 0
$ws.Range("D95").Value = 6741
$ws.Range("E95").Value = 1146
$ws.Range("F95").Value = 0
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = 63

$ws.Range("A96").Value = "Haiti"
$ws.Range("B96").Value = 7810
$ws.Range("C96").Value = 29
$ws.Range("D96").Value = 5123
$ws.Range("E96").Value = 2495
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 192

$ws.Range("A122").Value = "Sri Lanka"
$ws.Range("B122").Value = 2886
$ws.Range("C122").Value = 4
$ws.Range("D122").Value = 2658
$ws.Range("E122").Value = 217
$ws.Range("F122").Value = 0
$ws.Range("G122").Value = 0
$ws.Range("H122").Value = 11

$ws.Range("A135").Value = "Sierra Leona"
$ws.Range("B135").Value = 1947
$ws.Range("C135").Value = 7
$ws.Range("D135").Value = 1502
$ws.Range("E135").Value = 376
$ws.Range("F135").Value = 0
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 69

$ws.Range("A156").Value = "Principado de Andorra"
$ws.Range("B156").Value = 989
$ws.Range("C156").Value = 8
$ws.Range("D156").Value = 863
$ws.Range("E156").Value = 73
$ws.Range("F156").Value = 0
$ws.Range("G156").Value = 0
$ws.Range("H156").Value = 53

$ws.Range("A169").Value = "Burundi"
$ws.Range("B169").Value = 412
$ws.Range("C169").Value = 2
$ws.Range("D169").Value = 315
$ws.Range("E169").Value = 96
$ws.Range("F169").Value = 0
$ws.Range("G169").Value = 0
$ws.Range("H169").Value = 1

$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("B213").Value = 13
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 13
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Montserrat"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 12
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 1
